# Actualiza base de datos EC: intercambia los valores de "Periodo Mora" y
# "Valor Mora" entre las filas 16 y 17 (el periodo 1707 pasa a la fila 16 y
# el periodo 1708 pasa a la fila 17, llevando consigo su respectivo valor).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values before overwriting anything.
$periodo16 = $ws.Range("E16").Value()
$valor16   = $ws.Range("F16").Value()
$periodo17 = $ws.Range("E17").Value()
$valor17   = $ws.Range("F17").Value()

# Swap the "Periodo Mora" / "Valor Mora" pairs between the two rows.
$ws.Range("E16").Value = $periodo17
$ws.Range("F16").Value = $valor17
$ws.Range("E17").Value = $periodo16
$ws.Range("F17").Value = $valor16
